$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: headers (lowercase now), add D1 = "username" ---
$ws.Cells.Item(1,1).Value = "email"
$ws.Cells.Item(1,2).Value = "password"
$ws.Cells.Item(1,3).Value = "role"
$ws.Cells.Item(1,4).Value = "username"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").HorizontalAlignment = -4108

# --- Row 2: abc.123@gmail.com entry, add D2 = "ABC" ---
$ws.Cells.Item(2,4).Value = "ABC"

# --- Row 3: replace admin/admin/admin row with admin@gmail.com entry ---
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:admin@gmail.com")
$ws.Cells.Item(3,1).Value = "admin@gmail.com"
$ws.Cells.Item(3,2).Value = "admin"
$ws.Cells.Item(3,3).Value = "admin"
$ws.Cells.Item(3,4).Value = "Admin"
$ws.Range("A3").Style = "Hyperlink"

# --- Row 4: def.456@gmail.com entry, add D4 = "DEF" ---
$ws.Cells.Item(4,4).Value = "DEF"

# Update selection to match target (D2 selected)
$ws.Range("D2").Select()

$wb.Save()
